$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text formatting
# (values such as "2.20", "145.50", "0.100" must not be coerced to numbers).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.867.77'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.15%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.671.59'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.55%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.29%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.09'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.60%  '

# Row 6
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +5.55%  '

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.31%  '

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.34%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0619'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.92%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.33'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +3.98%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0893'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +4.05%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.907.78'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.88%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.694.80'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.44%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.09'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.68%  '

# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.28%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.66'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.50%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.896.98'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.09%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '232.59'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -3.99%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.83'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.57%  '

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.02%  '

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.17%  '

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.66%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.19'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.55%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.20'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -3.24%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.50'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.58%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.116'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.39%  '

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.46%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.96'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.67%  '

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.16%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0498'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.36%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.18'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.90%  '

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.55%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.472.04'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.98%  '

# Row 34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +4.23%  '

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +2.73%  '

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.15%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.901'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.90%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.571'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.90%  '

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.39%  '

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.30%  '

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.18%  '

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +3.61%  '

# Row 43
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '65.80'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.21%  '

# Row 44
$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.971'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +6.77%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.816.61'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.97%  '

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.18%  '

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.03%  '

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.46%  '

# Row 49
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.100'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.51%  '

# Row 50
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0508'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.22%  '

# Row 51
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.58'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.98%  '
